$p = $ppt.ActivePresentation

# --- Slide 4 (repo path) == physical slide12.xml: "Solutions trouvées" ---
$s4 = $p.Slides.Item(12)

# Remove the three small diagram pictures that used to sit under "Partie 1".
foreach ($n in @("Image 3", "Image 5", "Image 26")) {
    $s4.Shapes.Item($n).Delete()
}

# Reposition + retitle the "Partie 1" caption textbox.
$caption = $s4.Shapes.Item("ZoneTexte 6")
$caption.Left = 4115087 / 12700
$caption.Top = 1292497 / 12700
$caption.TextFrame.TextRange.Text = "Partie 1 : Pluviomètre"

# --- Slide 5 (repo path): acquisition slide, now hidden & stripped down ---
$s5 = $p.Slides.Item(5)

# Hide the slide from the slide show (adds show="0").
$s5.SlideShowTransition.Hidden = $true

# Drop the click-animation effects tied to the shapes we are about to remove,
# otherwise the <p:timing>/<p:bldLst> trees would be left dangling.
$seq = $s5.TimeLine.MainSequence
for ($i = $seq.Count; $i -ge 1; $i--) {
    $seq.Item($i).Delete()
}

# Remove the leftover acquisition-workflow pictures and text boxes.
$namesToRemove = @(
    "Image 28",
    "Image 24",
    "Image 34",
    "Image 3",
    "ZoneTexte 1",
    "ZoneTexte 4",
    "ZoneTexte 5",
    "ZoneTexte 6",
    "ZoneTexte 10",
    "ZoneTexte 11",
    "Image 13"
)
foreach ($n in $namesToRemove) {
    $s5.Shapes.Item($n).Delete()
}
